# Fixed Grogu information and added dynamic font size.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "First Screen Appearance" text for Grogu (row 6, column F) -
# correct the punctuation from "(Season 1, 2019)" to ": Season 1 (2019)"
$ws.Range("F6").Value = "The Mandalorian: Season 1 (2019)"

# Move the current selection (as last left by the editing session) to F9
$ws.Activate()
$ws.Range("F9").Select()
